$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.149.31"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "2.313.05"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'310.09"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  +6.32%  "
$ws.Range("D7").Value = "'0.537"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +7.85%  "
$ws.Range("D10").Value = "'36.05"
$ws.Range("E10").Value = "  +4.16%  "
$ws.Range("D11").Value = "'0.0816"
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +6.65%  "
$ws.Range("D14").Value = "2.668.13"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").Value = "'14.99"
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("D16").Value = "2.305.55"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "'0.813"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "43.056.04"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "'68.51"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'241.25"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").Value = "'2.03"
$ws.Range("E24").Value = "  +6.36%  "
$ws.Range("E25").Value = "  +3.65%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'24.85"
$ws.Range("E27").Value = "  +5.70%  "
$ws.Range("D28").Value = "'37.37"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("D29").Value = "'9.69"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "'166.76"
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'3.16"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +6.82%  "
$ws.Range("D36").Value = "'0.0744"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("E37").Value = "  +3.48%  "
$ws.Range("D38").Value = "'2.39"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").Value = "'19.53"
$ws.Range("E43").Value = "  +5.10%  "
$ws.Range("D44").Value = "'0.0291"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "1.978.85"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "'2.98"
$ws.Range("E48").Value = "  +19.12%  "
$ws.Range("D49").Value = "'55.66"
$ws.Range("E49").Value = "  +6.05%  "
$ws.Range("D50").Value = "2.539.80"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("E51").Value = "  +4.21%  "
